$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 6478
$ws.Range("D25").Value = 6028208
$ws.Range("E25").Value = 930.5662241432541
$ws.Range("F25").Value = 9.964352401969112
$ws.Range("H25").Value = 26.24926542133505
